# Apply updated Betfair back/lay odds values for 2026-02-05 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("S2").Value = 2.62
$ws.Range("Z2").Value = 970

# Row 3
$ws.Range("P3").Value = 1.4
$ws.Range("Q3").Value = 3.05
$ws.Range("V3").Value = 1.38

# Row 4
$ws.Range("L4").Value = 1.29
$ws.Range("S4").Value = 2.82
$ws.Range("T4").Value = 1.63
$ws.Range("U4").Value = 1.73

# Row 6
$ws.Range("Q6").Value = 1.38
$ws.Range("T6").Value = 1.95

# Row 7
$ws.Range("Q7").Value = 1.36

# Row 10
$ws.Range("F10").Value = 1.72
$ws.Range("G10").Value = 1.91
$ws.Range("K10").Value = 4
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 2.52
$ws.Range("P10").Value = 1.71
$ws.Range("Q10").Value = 2.12
$ws.Range("R10").Value = 1.18
$ws.Range("S10").Value = 3.55
$ws.Range("T10").Value = 1.05
$ws.Range("U10").Value = 1.7
$ws.Range("W10").Value = 2.1
$ws.Range("X10").Value = 15

# Row 11
$ws.Range("J11").Value = 2.82
$ws.Range("M11").Value = 1.13
$ws.Range("S11").Value = 5.6

# Row 12
$ws.Range("H12").Value = 3.7
$ws.Range("K12").Value = 3.6
$ws.Range("L12").Value = 1.45
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 3.2
$ws.Range("O12").Value = 1.35
$ws.Range("R12").Value = 1.28
$ws.Range("S12").Value = 3.25
$ws.Range("T12").Value = 1.05
$ws.Range("U12").Value = 1.9
$ws.Range("V12").Value = 1.35
$ws.Range("W12").Value = 1.76
$ws.Range("X12").Value = 15.5
$ws.Range("Y12").Value = 15.5
$ws.Range("Z12").Value = 38
$ws.Range("AA12").Value = 90
$ws.Range("AB12").Value = 10
$ws.Range("AC12").Value = 8.199999999999999
$ws.Range("AD12").Value = 18.5
$ws.Range("AE12").Value = 55
$ws.Range("AF12").Value = 17.5
$ws.Range("AG12").Value = 15.5
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 70
$ws.Range("AJ12").Value = 36
$ws.Range("AK12").Value = 30
$ws.Range("AL12").Value = 50
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 25
$ws.Range("AO12").Value = 1000

# Row 13
$ws.Range("F13").Value = 2.12
$ws.Range("G13").Value = 2.38
$ws.Range("H13").Value = 3.7
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 3.05
$ws.Range("K13").Value = 3.55
$ws.Range("L13").Value = 1.01
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 2.38
$ws.Range("O13").Value = 1.01
$ws.Range("P13").Value = 1.62
$ws.Range("Q13").Value = 2.3
$ws.Range("R13").Value = 1.16
$ws.Range("S13").Value = 3.45
$ws.Range("T13").Value = 1.67
$ws.Range("U13").Value = 1.6
$ws.Range("V13").Value = 1.25
$ws.Range("W13").Value = 1.72
$ws.Range("X13").Value = 15
$ws.Range("Y13").Value = 17.5
$ws.Range("Z13").Value = 40
$ws.Range("AA13").Value = 1000
$ws.Range("AB13").Value = 11
$ws.Range("AC13").Value = 10.5
$ws.Range("AD13").Value = 25
$ws.Range("AE13").Value = 90
$ws.Range("AF13").Value = 19
$ws.Range("AG13").Value = 16
$ws.Range("AH13").Value = 30
$ws.Range("AI13").Value = 100
$ws.Range("AJ13").Value = 46
$ws.Range("AK13").Value = 40
$ws.Range("AL13").Value = 70
$ws.Range("AM13").Value = 1000
$ws.Range("AN13").Value = 38
$ws.Range("AO13").Value = 1000

# Row 14
$ws.Range("F14").Value = 1.67
$ws.Range("G14").Value = 1.7
$ws.Range("J14").Value = 3.8
$ws.Range("K14").Value = 4.1
$ws.Range("L14").Value = 1.43
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 3.4
$ws.Range("O14").Value = 1.34
$ws.Range("R14").Value = 1.33
$ws.Range("S14").Value = 3.6
$ws.Range("T14").Value = 1.01
$ws.Range("U14").Value = 1.75
$ws.Range("V14").Value = 1.17
$ws.Range("W14").Value = 2.42
$ws.Range("X14").Value = 13.5
$ws.Range("Y14").Value = 20
$ws.Range("Z14").Value = 55
$ws.Range("AA14").Value = 220
$ws.Range("AB14").Value = 7.8
$ws.Range("AC14").Value = 8.800000000000001
$ws.Range("AD14").Value = 24
$ws.Range("AE14").Value = 110
$ws.Range("AF14").Value = 10.5
$ws.Range("AG14").Value = 9.800000000000001
$ws.Range("AH14").Value = 24
$ws.Range("AI14").Value = 110
$ws.Range("AJ14").Value = 16
$ws.Range("AK14").Value = 18.5
$ws.Range("AL14").Value = 42
$ws.Range("AM14").Value = 160
$ws.Range("AN14").Value = 12.5
$ws.Range("AO14").Value = 150
